# Insert a new weekly record row at row 75 on "Sheet1", shifting the
# existing rows (old 75-89) down to (76-90), and fill in the new row's
# data: Poroto granado, Región del Maule, fecha 44924 (2022-12-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 75; existing data shifts
# down automatically (old row 75 -> 76, ..., old row 89 -> 90).
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record.
$ws.Cells.Item(75, 1).Value = 10
$ws.Cells.Item(75, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(75, 3).Value = 'La Araucanía'
$ws.Cells.Item(75, 4).Value = 44924
$ws.Cells.Item(75, 5).Value = 9
$ws.Cells.Item(75, 6).Value = 100112030
$ws.Cells.Item(75, 7).Value = 'Poroto granado'
$ws.Cells.Item(75, 8).Value = 'Sin especificar'
$ws.Cells.Item(75, 9).Value = 'Primera'
$ws.Cells.Item(75, 10).Value = 205
$ws.Cells.Item(75, 11).Value = 43000
$ws.Cells.Item(75, 12).Value = 45000
$ws.Cells.Item(75, 13).Value = 44220
$ws.Cells.Item(75, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(75, 15).Value = 'Región del Maule'
$ws.Cells.Item(75, 16).Value = 1769
$ws.Cells.Item(75, 17).Value = 25
$ws.Cells.Item(75, 18).Value = 'Hortaliza'
